# 11/18/2020 Experiment results were updated.
# Adds a new "Turbidity" worksheet (containing the newly-collected
# experiment data) at the end of the workbook, and makes it the active sheet.

$wb = $excel.ActiveWorkbook

# New turbidity measurements (9 columns x 12 rows).
$data = @(
  @(0.0179,0.0184,0.0137,0.0107,0.0108,0.0132,0.0193,0.0194,0.015),
  @(0.013,0.0166,0.0086,0.0087,0.0091,0.0112,0.0169,0.0173,0.0132),
  @(0.0142,0.041,0.0047,0.005,0.0023,0.0051,0.0408,0.0409,0.0077),
  @(0.0167,0.027,0.0112,0.0083,0.0096,0.0117,0.0276,0.028,0.0165),
  @(0.0165,0.049,0.0064,0.0032,0.0058,0.0093,0.0482,0.0489,0.0188),
  @(0.0142,0.0131,0.0059,0.0057,0.0054,0.0107,0.0142,0.0141,0.0138),
  @(0.0098,0.0132,0.0022,0,0.0009,0.0073,0.0138,0.0141,0.0096),
  @(0.014,0.021,0.0071,0.0042,0.0106,0.0134,0.0215,0.0218,0.0176),
  @(0.0329,0.0159,0.0355,0.0325,0.0199,0.0246,0.0164,0.0165,0.032),
  @(0.0043,0.0145,0.008,0.0058,0.0023,0.005,0.0148,0.015,0.0043),
  @(0.0602,0.0334,0.0667,0.0725,0.0435,0.0406,0.0333,0.0337,0.0435),
  @(0.0173,0.0169,0.0195,0.018,0.0081,0.01,0.0174,0.0176,0.0207)
)

# Add the new sheet after the last existing sheet (so it lands after "OV").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Turbidity"

for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt $data[$r].Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
    }
}

# Matches the saved selection state for the new sheet.
$ws.Range("E16").Select()
